# Append the new portfolio data row (row 61) to Sheet1, matching the
# plain/unstyled format used by the existing data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$row = 61

# Force the date column to be stored as literal text (not auto-converted
# to a date serial number) while keeping the same "no explicit style"
# appearance as the other data rows.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-10-15"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 53.72999954223633
$ws.Cells.Item($row, 3).Value = 390.8500061035156
$ws.Cells.Item($row, 4).Value = 354.3500061035156
